$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATA_POINTS (row 3) and SAMPLE_POINTS (row 4) from 100 to 5
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 5

# Update the PROMPT text (row 8) - remove "NOT paraphrase or change any part of the text except for"
$ws.Range("B8").Value = "You are an advanced anonymizer that replaces personally identifiable information (PII) with a category label. You will replacing PII with its category in square brackets.`r`n`r`nExample:`r`nInput: My name is Alice and I live in London.`r`nOutput: My name is [NAME] and I live in [LOCATION]."
